# Adds two new species-observation rows (72 and 73) to the "Artfynd" sheet,
# matching the appended rows from the source diff. Extends the used range
# from A1:AY71 to A1:AY73.
#
# Note: a few string values look like numbers/dates ("1", "2023-09-14").
# Plain Excel COM type-inference would silently turn those into a real
# number / date serial on assignment, which would not match the source
# data (stored as plain text). Prefixing with a leading apostrophe forces
# Excel to keep them as text, exactly like typing '2023-09-14 into a cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 72 : Craterellus lutescens (Rodgul trumpetsvamp) ---
$ws.Range("A72").Value  = 112090588
$ws.Range("B72").Value  = 89183
$ws.Range("C72").Value  = 'Ovaliderad'
$ws.Range("D72").Value  = 'LC'
$ws.Range("E72").Value  = 3215
$ws.Range("F72").Value  = 'Rödgul trumpetsvamp'
$ws.Range("G72").Value  = 'Craterellus lutescens'
$ws.Range("H72").Value  = '(Fr.) Fr.'
$ws.Range("P72").Value  = 'Knivsta (Knivsta), Upl'
$ws.Range("Q72").Value  = 654788.2679259261
$ws.Range("R72").Value  = 6626333.524893245
$ws.Range("S72").Value  = 20
$ws.Range("T72").Value  = 'Uppsala'
$ws.Range("U72").Value  = 'Knivsta'
$ws.Range("V72").Value  = 'Uppland'
$ws.Range("W72").Value  = 'Alsike'
$ws.Range("Y72").Value  = "'2023-09-14"
$ws.Range("Z72").Value  = '15:59'
$ws.Range("AA72").Value = "'2023-09-14"
$ws.Range("AB72").Value = '15:59'
$ws.Range("AD72").Value = $false
$ws.Range("AE72").Value = $false
$ws.Range("AG72").Value = $false
$ws.Range("AW72").Value = 'Marie Kvarnström'
$ws.Range("AX72").Value = 'Marie Kvarnström'

# --- row 73 : Actaea spicata (Svart trolldruva) ---
$ws.Range("A73").Value  = 112090750
$ws.Range("B73").Value  = 98446
$ws.Range("C73").Value  = 'Ovaliderad'
$ws.Range("D73").Value  = 'LC'
$ws.Range("E73").Value  = 222771
$ws.Range("F73").Value  = 'Svart trolldruva'
$ws.Range("G73").Value  = 'Actaea spicata'
$ws.Range("H73").Value  = 'L.'
$ws.Range("I73").Value  = "'1"
$ws.Range("J73").Value  = 'plantor/tuvor'
$ws.Range("K73").Value  = 'i frukt'
$ws.Range("P73").Value  = 'Knivsta (Knivsta), Upl'
$ws.Range("Q73").Value  = 654798.0062938032
$ws.Range("R73").Value  = 6626354.618875842
$ws.Range("S73").Value  = 20
$ws.Range("T73").Value  = 'Uppsala'
$ws.Range("U73").Value  = 'Knivsta'
$ws.Range("V73").Value  = 'Uppland'
$ws.Range("W73").Value  = 'Alsike'
$ws.Range("Y73").Value  = "'2023-09-14"
$ws.Range("Z73").Value  = '16:02'
$ws.Range("AA73").Value = "'2023-09-14"
$ws.Range("AB73").Value = '16:02'
$ws.Range("AD73").Value = $false
$ws.Range("AE73").Value = $false
$ws.Range("AG73").Value = $false
$ws.Range("AW73").Value = 'Marie Kvarnström'
$ws.Range("AX73").Value = 'Marie Kvarnström'
